$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "27.121.39"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "  +0.43%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "1.827.05"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "  +0.20%  "
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'1.008"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "  +0.26%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'312.46"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "  +0.35%  "
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "  +0.21%  "
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'0.4696"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "  +0.09%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'0.3665"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "  +0.08%  "
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.07389"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "  +0.45%  "
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'0.8797"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "  +0.57%  "
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'20.29"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "  -0.06%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "1.900.33"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "  +4.02%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'0.07711"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "  +5.63%  "
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = "  +1.69%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'5.370"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "  -1.10%  "
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'6.530"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "  +0.17%  "
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "  +0.03%  "
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'0.000008719"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "  -0.30%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'1.007"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "  +0.34%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "27.643.33"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "  +2.32%  "
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'14.61"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "  -0.57%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'5.240"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "  -0.97%  "
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'10.62"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "  -0.01%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "2.085.73"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "  +1.52%  "
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'1.872"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "  -1.17%  "
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "  -0.17%  "
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'18.48"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "  +0.30%  "
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.Value = "'2.137"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "  -0.21%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = "'5.180"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "  -1.27%  "
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = "'116.50"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "  -0.30%  "
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = "'0.08922"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "  +0.37%  "
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = "'0.7444"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "  -1.41%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = "'1.163"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "  +0.16%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = "'4.514"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "  +0.09%  "
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.Value = "'2.942"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "  +0.37%  "
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.Value = "'2.635"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "  +11.15%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'1.008"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "  +0.29%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'1.090"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "  -0.65%  "
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.05302"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "  -0.25%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.01934"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "  -0.81%  "
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'7.334"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "  +1.79%  "
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'2.927"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "  -1.90%  "
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'0.5264"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "  -0.76%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'0.1643"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "  -0.68%  "
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'8.392"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "  -1.03%  "
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'0.4908"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "  +0.24%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'10.39"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'1.007"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "  +0.30%  "
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'104.50"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "  +1.20%  "
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'1.653"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "  -0.58%  "
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'0.06277"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "  -0.39%  "
$c.Style = "Normal"
